# "Generate Report for Handback" - refresh the handback-status timestamps.
#
# The CI handback report recorded new handoff/handback timestamps for the
# 9b365740-...-md entry (the c0bdec6a-...-md entry's timestamps did not
# change in this run):
#   - Overview!G2            "Latest HO Xliff Generate Date"      -> 2016-09-06 17:25:20
#   - zh-cn!H2  (row for 9b365740...) "Correspond Handoff Datetime"  -> 2016-09-06 17:25:02
#   - zh-cn!K2  (row for 9b365740...) "Correspond Handback DateTime" -> 2016-09-06 17:25:37
#   - de-de!K2  (row for 9b365740...) "Correspond Handback DateTime" -> 2016-09-06 17:25:45

$wb = $excel.ActiveWorkbook
$dtFormat = "yyyy-mm-dd HH:mm:ss"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-09-06 17:25:20"
$overview.Range("G2").NumberFormat = $dtFormat

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-09-06 17:25:02"
$zhcn.Range("H2").NumberFormat = $dtFormat
$zhcn.Range("K2").Value = "2016-09-06 17:25:37"
$zhcn.Range("K2").NumberFormat = $dtFormat

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("K2").Value = "2016-09-06 17:25:45"
$dede.Range("K2").NumberFormat = $dtFormat
